# SewerageTaxTestData.xlsx edit
# - Rename header C1 to "noOfClosetsResidential"
# - Insert new column D "noOfClosetsNonResidential" (shifts old D "documentNumber" to E)
# - For RESIDENTIAL rows (2,3): keep C numeric, set D = "null" (text)
# - For NON RESIDENTIAL rows (4,5): set C = "null" (text), D = numeric (old C values)
# - Add new rows 6,7 for MIXED property type with numeric-but-text-formatted C/D values
# - Adjust column widths

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before D so documentNumber shifts from D to E
$ws.Range("D1").EntireColumn.Insert()

# 2. Update headers
$ws.Range("C1").Value = "noOfClosetsResidential"
$ws.Range("D1").Value = "noOfClosetsNonResidential"

# 3. Fix up rows 2-3 (RESIDENTIAL): D should be "null" string (C already numeric & correct)
$ws.Range("D2").Value = "null"
$ws.Range("D3").Value = "null"

# 4. Fix up rows 4-5 (NON RESIDENTIAL): move numeric closet count from C to D, set C = "null"
$c4 = $ws.Range("C4").Value()
$c5 = $ws.Range("C5").Value()
$ws.Range("D4").Value = $c4
$ws.Range("D5").Value = $c5
$ws.Range("C4").Value = "null"
$ws.Range("C5").Value = "null"

# 5. Add new rows for MIXED property type
$ws.Range("A6").Value = "creation3"
$ws.Range("B6").Value = "MIXED"
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 123

$ws.Range("A7").Value = "change3"
$ws.Range("B7").Value = "MIXED"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 456

# Apply text number format to the MIXED rows' closet count cells
$ws.Range("C6:D7").NumberFormat = "@"

# 6. Adjust column widths to match target layout (nearest achievable values;
#    the host engine rounds ColumnWidth to whole pixels internally)
$ws.Columns.Item(1).ColumnWidth = 17.0
$ws.Columns.Item(2).ColumnWidth = 18.0
$ws.Columns.Item(3).ColumnWidth = 16.0
$ws.Columns.Item(4).ColumnWidth = 31.3333333333333
$ws.Columns.Item(5).ColumnWidth = 21.3333333333333
$ws.Columns.Item(6).ColumnWidth = 7.66666666666667

# 7. Update selection to D7 like the target file
$ws.Range("D7").Select()
